# Update cryptocurrency price/volume data per latest scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking Price values to stay text (matches source formatting,
# e.g. "1.000" / "7.440" must not collapse to 1 / 7.44)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "29.312.15"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.840.94"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "239.45"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "0.6279"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.07517"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "24.46"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "0.07684"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "1.859.64"
$ws.Range("E12").Value = "  -6.32%  "
$ws.Range("D13").Value = "4.968"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "0.6768"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "0.00001023"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").Value = "82.95"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "2.101.22"
$ws.Range("E17").Value = "  -7.20%  "
$ws.Range("D18").Value = "6.124"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "29.343.20"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "228.09"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "7.440"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "156.76"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "0.1389"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "8.342"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").Value = "1.458"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "1.268"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "0.05638"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").Value = "4.107"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "4.024"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("D36").Value = "0.7140"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "1.240.97"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "0.01804"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "2.765"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "6.227"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").Value = "0.9018"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "101.60"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "65.40"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").Value = "0.00000000118"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").Value = "7.072"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").Value = "0.3987"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.673"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.900"
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").Value = "0.1117"
$ws.Range("E51").Value = "  -0.22%  "
